$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "INV001"
$ws.Range("B2").Value = "Inventaire du 2025-06-02_1131"
$ws.Range("C2").Value = "2025-06-02 11:31:22"
$ws.Range("D2").Value = "En préparation"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = "Utilisateur"
